# Update cryptocurrency price ("D") and 1h volume change ("E") columns
# to reflect the refreshed values from the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number (e.g. "300.16") are written
# with a leading apostrophe so Excel keeps them as text instead of silently
# converting them to a floating point Number (which would corrupt values such
# as "43.033.52" or drop trailing/insignificant digits like "300.16").
$updates = @(
    @{ Cell = 'D2'; Value = '43.033.52' }
    @{ Cell = 'E2'; Value = '  +0.00%  ' }
    @{ Cell = 'D3'; Value = '2.300.03' }
    @{ Cell = 'E3'; Value = '  +0.04%  ' }
    @{ Cell = 'E4'; Value = '  +0.01%  ' }
    @{ Cell = 'D5'; Value = '''300.16' }
    @{ Cell = 'E5'; Value = '  -0.12%  ' }
    @{ Cell = 'D6'; Value = '''97.82' }
    @{ Cell = 'E6'; Value = '  -1.65%  ' }
    @{ Cell = 'E7'; Value = '  +2.63%  ' }
    @{ Cell = 'E8'; Value = '  -0.04%  ' }
    @{ Cell = 'E9'; Value = '  +0.79%  ' }
    @{ Cell = 'D10'; Value = '''36.06' }
    @{ Cell = 'E10'; Value = '  -0.44%  ' }
    @{ Cell = 'E11'; Value = '  -0.03%  ' }
    @{ Cell = 'E12'; Value = '  +0.78%  ' }
    @{ Cell = 'D13'; Value = '''17.71' }
    @{ Cell = 'E13'; Value = '  -0.94%  ' }
    @{ Cell = 'D14'; Value = '''6.88' }
    @{ Cell = 'E14'; Value = '  -0.24%  ' }
    @{ Cell = 'D15'; Value = '2.657.77' }
    @{ Cell = 'D16'; Value = '2.255.28' }
    @{ Cell = 'E16'; Value = '  -1.84%  ' }
    @{ Cell = 'D17'; Value = '''0.788' }
    @{ Cell = 'E17'; Value = '  -1.39%  ' }
    @{ Cell = 'D18'; Value = '42.929.80' }
    @{ Cell = 'E18'; Value = '  -0.05%  ' }
    @{ Cell = 'D19'; Value = '''12.76' }
    @{ Cell = 'E19'; Value = '  -0.57%  ' }
    @{ Cell = 'D20'; Value = '0.0₃0909' }
    @{ Cell = 'E20'; Value = '  +0.47%  ' }
    @{ Cell = 'E21'; Value = '  +0.32%  ' }
    @{ Cell = 'D22'; Value = '''68.79' }
    @{ Cell = 'E22'; Value = '  +1.28%  ' }
    @{ Cell = 'D23'; Value = '''237.82' }
    @{ Cell = 'E23'; Value = '  +0.78%  ' }
    @{ Cell = 'E24'; Value = '  -1.09%  ' }
    @{ Cell = 'E25'; Value = '  -0.40%  ' }
    @{ Cell = 'E26'; Value = '  -0.45%  ' }
    @{ Cell = 'D27'; Value = '''4.01' }
    @{ Cell = 'E27'; Value = '  -0.18%  ' }
    @{ Cell = 'E28'; Value = '  +0.26%  ' }
    @{ Cell = 'D29'; Value = '''164.72' }
    @{ Cell = 'E29'; Value = '  -2.61%  ' }
    @{ Cell = 'D30'; Value = '''2.04' }
    @{ Cell = 'E30'; Value = '  -0.29%  ' }
    @{ Cell = 'D31'; Value = '''9.12' }
    @{ Cell = 'E31'; Value = '  -0.10%  ' }
    @{ Cell = 'D32'; Value = '''33.01' }
    @{ Cell = 'E32'; Value = '  -4.03%  ' }
    @{ Cell = 'E33'; Value = '  +0.00%  ' }
    @{ Cell = 'E34'; Value = '  +1.38%  ' }
    @{ Cell = 'E35'; Value = '  +4.18%  ' }
    @{ Cell = 'D36'; Value = '''17.89' }
    @{ Cell = 'E36'; Value = '  +1.52%  ' }
    @{ Cell = 'E37'; Value = '  -0.11%  ' }
    @{ Cell = 'E38'; Value = '  +1.32%  ' }
    @{ Cell = 'E39'; Value = '  +0.38%  ' }
    @{ Cell = 'E40'; Value = '  -0.43%  ' }
    @{ Cell = 'E41'; Value = '  -0.93%  ' }
    @{ Cell = 'E42'; Value = '  +1.02%  ' }
    @{ Cell = 'D43'; Value = '2.016.69' }
    @{ Cell = 'E43'; Value = '  +1.61%  ' }
    @{ Cell = 'E44'; Value = '  -1.79%  ' }
    @{ Cell = 'D45'; Value = '''2.19' }
    @{ Cell = 'E45'; Value = '  -4.84%  ' }
    @{ Cell = 'D46'; Value = '''10.39' }
    @{ Cell = 'E46'; Value = '  +2.22%  ' }
    @{ Cell = 'D47'; Value = '''17.52' }
    @{ Cell = 'E47'; Value = '  -0.46%  ' }
    @{ Cell = 'E48'; Value = '  -2.20%  ' }
    @{ Cell = 'D49'; Value = '''54.11' }
    @{ Cell = 'E49'; Value = '  -2.73%  ' }
    @{ Cell = 'D50'; Value = '2.526.81' }
    @{ Cell = 'E50'; Value = '  +0.10%  ' }
    @{ Cell = 'E51'; Value = '  -1.28%  ' }
)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}
